# "Cleaning ver data sheets"
#
# The sheet originally laid out three (k, m) column-pairs side by side with
# a narrow blank spacer column between each pair:
#   A,B = pair1 (k,m) | C = spacer | D,E = pair2 (k,m) | F = spacer | G,H = pair3 (k,m)
#
# The cleaned-up layout removes the two blank spacer columns so the three
# pairs sit directly next to one another:
#   A,B = pair1 (k,m) | C,D = pair2 (k,m) | E,F = pair3 (k,m)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first spacer column (originally column C). Everything to the
# right (the second pair + second spacer + third pair) shifts one column left.
$ws.Columns("C").Delete()

# The second spacer column was originally column F; after the first delete
# it is now column E. Remove it too, shifting the third pair left again.
$ws.Columns("E").Delete()

# Update the active selection to match the cleaned-up sheet.
$ws.Range("J8").Select()
